$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '68.970.35'
$ws.Range('E2').Value = '  -1.94%  '

# Row 3
$ws.Range('D3').Value = '3.520.36'
$ws.Range('E3').Value = '  -2.91%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = '588.60'
$ws.Range('E5').Value = '  +1.53%  '

# Row 6
$ws.Range('D6').Value = '170.54'
$ws.Range('E6').Value = '  -2.47%  '

# Row 7
$ws.Range('E7').Value = '  +0.63%  '

# Row 8
$ws.Range('D8').Value = '3.513.63'
$ws.Range('E8').Value = '  -2.83%  '

# Row 9
$ws.Range('E9').Value = '  +0.00%  '

# Row 10
$ws.Range('E10').Value = '  -3.54%  '

# Row 11
$ws.Range('D11').Value = '6.84'
$ws.Range('E11').Value = '  -1.70%  '

# Row 12
$ws.Range('E12').Value = '  -4.02%  '

# Row 13
$ws.Range('D13').Value = '46.97'
$ws.Range('E13').Value = '  -2.30%  '

# Row 14
$ws.Range('E14').Value = '  -2.21%  '

# Row 15
$ws.Range('D15').Value = '4.078.94'
$ws.Range('E15').Value = '  -3.10%  '

# Row 16
$ws.Range('E16').Value = '  -4.70%  '

# Row 17
$ws.Range('D17').Value = '616.03'
$ws.Range('E17').Value = '  -8.57%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.523.39'
$ws.Range('E18').Value = '  -2.99%  '

# Row 19
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '69.060.47'
$ws.Range('E19').Value = '  -1.93%  '

# Row 20
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('D21').Value = '17.36'
$ws.Range('E21').Value = '  -2.12%  '

# Row 22
$ws.Range('D22').Value = '11.11'
$ws.Range('E22').Value = '  -2.45%  '

# Row 23
$ws.Range('E23').Value = '  -5.69%  '

# Row 24
$ws.Range('D24').Value = '15.77'
$ws.Range('E24').Value = '  -7.29%  '

# Row 25
$ws.Range('D25').Value = '96.53'
$ws.Range('E25').Value = '  -3.07%  '

# Row 26
$ws.Range('E26').Value = '  -1.38%  '

# Row 28
$ws.Range('E28').Value = '  -5.80%  '

# Row 29
$ws.Range('E29').Value = '  -6.33%  '

# Row 30
$ws.Range('E30').Value = '  -5.35%  '

# Row 31
$ws.Range('D31').Value = '3.12'
$ws.Range('E31').Value = '  -5.34%  '

# Row 32
$ws.Range('D32').Value = '8.48'
$ws.Range('E32').Value = '  -5.55%  '

# Row 33
$ws.Range('E33').Value = '  -4.88%  '

# Row 34
$ws.Range('E34').Value = '  -7.85%  '

# Row 35
$ws.Range('D35').Value = '626.36'
$ws.Range('E35').Value = '  +8.37%  '

# Row 36
$ws.Range('D36').Value = '10.73'
$ws.Range('E36').Value = '  -2.61%  '

# Row 37
$ws.Range('D37').Value = '3.46'
$ws.Range('E37').Value = '  -12.07%  '

# Row 38
$ws.Range('E38').Value = '  -4.09%  '

# Row 39
$ws.Range('D39').Value = '57.31'
$ws.Range('E39').Value = '  -1.38%  '

# Row 40
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.04%  '

# Row 41
$ws.Range('D41').Value = '0.0446'
$ws.Range('E41').Value = '  -1.05%  '

# Row 42
$ws.Range('E42').Value = '  -3.98%  '

# Row 43
$ws.Range('D43').Value = '3.376.70'
$ws.Range('E43').Value = '  -4.91%  '

# Row 44
$ws.Range('E44').Value = '  -4.84%  '

# Row 45
$ws.Range('D45').Value = '32.69'
$ws.Range('E45').Value = '  -5.00%  '

# Row 46
$ws.Range('D46').Value = '0.0₃0694'
$ws.Range('E46').Value = '  -4.71%  '

# Row 47
$ws.Range('E47').Value = '  -5.25%  '

# Row 48
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').Value = '  -1.77%  '

# Row 49
$ws.Range('E49').Value = '  -2.52%  '

# Row 50
$ws.Range('D50').Value = '133.53'
$ws.Range('E50').Value = '  -2.07%  '

# Row 51
$ws.Range('E51').Value = '  +12.62%  '
